$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C entirely (it is removed in the target state)
$ws.Columns("C:C").Delete() | Out-Null

# Update column B values (rows 2-33); B1 stays unchanged at 1
$ws.Range("B2").Value = 0.07827832446500971
$ws.Range("B3").Value = 25.85178379958638
$ws.Range("B4").Value = 151.6999506794686
$ws.Range("B5").Value = 0.1857483474017147
$ws.Range("B6").Value = 2.259983897591462
$ws.Range("B7").Value = 3.45764265479874
$ws.Range("B8").Value = 28.34305588244797
$ws.Range("B9").Value = 0.996206034574341
$ws.Range("B10").Value = 1.01519036569497
$ws.Range("B11").Value = 2.148267168927628
$ws.Range("B12").Value = 0.9142165883506808
$ws.Range("B13").Value = 0.07829918587575628
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 0.5570218339441547
$ws.Range("B17").Value = 0.1789698279133011
$ws.Range("B18").Value = 0.1251479657877091
$ws.Range("B19").Value = 5.109587219701581
$ws.Range("B20").Value = 0.0117340989070336
$ws.Range("B21").Value = 0.4685710355772409
$ws.Range("B22").Value = 28.35066765640477
$ws.Range("B23").Value = 29.29956887469033
$ws.Range("B24").Value = 0.08161019259709537
$ws.Range("B25").Value = 0.4545683228568137
$ws.Range("B26").Value = 1.706235270895844
$ws.Range("B27").Value = 3.596233485751527
$ws.Range("B28").Value = 12.06697641810684
$ws.Range("B29").Value = 33.28598225363204
$ws.Range("B30").Value = 73817.19939611075
$ws.Range("B31").Value = 6.095572920455203
$ws.Range("B32").Value = 88.11583049130799
$ws.Range("B33").Value = 50.001183795832
